# Apply restated financial figures (Total Revenue through Change In Cash)
# to the DTEA yearly financials sheet, per the "Doing Updates for Financials" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    8 = @{ D=166700; E=160700; F=134500; G=105600; H=80500; I=54400 }
    9 = @{ D=86900; E=80000; F=63500; G=47800; H=36000; I=23900 }
    10 = @{ D=79800; E=80700; F=70900; G=57800; H=44500; I=30400 }
    14 = @{ D=11300; E=5900; G=3100 }
    15 = @{ D=7400; E=6600; F=4800; G=4100; H=3500 }
    17 = @{ D=185100; E=165400; F=123100; G=98300; H=75000; I=51700 }
    18 = @{ D=-18400; E=-4700; F=11300; G=7300; H=5500; I=2600 }
    20 = @{ F=-104900; H=-6300; I=-3300 }
    21 = @{ D=-10600; E=2200; F=-88700; G=10400; H=2700 }
    22 = @{ G=1700 }
    23 = @{ D=-19700; E=-4400; F=-94300; G=4600; H=-2300 }
    24 = @{ F=3500; G=-200; H=2300 }
    26 = @{ D=-19700; E=-2700; F=-97800; G=4800; H=-4600; I=-3200 }
    27 = @{ D=-19700; E=-2700; F=-97800; G=4800; H=-4600; I=-3200 }
    32 = @{ F=104900; H=6300; I=3300 }
    33 = @{ D=-21200; E=-2700; F=-97800; G=4800; H=-4600; I=-3200 }
    35 = @{ D=-21200; E=-2700; F=-97800; G=4800; H=-4600; I=-3200 }
    41 = @{ D=47200; E=48000; F=54000; G=14700; H=11400; I=5600 }
    43 = @{ D=4500; E=3000; G=2400; H=800 }
    44 = @{ D=18200; E=23300; F=13200; G=9300; H=8400; I=8000 }
    45 = @{ D=5700; E=4500; F=5900; H=2100; I=2000 }
    46 = @{ D=75700; E=78800; F=75600; G=28700; H=22800; I=16300 }
    48 = @{ D=27200; E=38100; F=35200; G=26500; H=22200; I=19300 }
    49 = @{ D=3300; E=2200; G=1200; H=1100 }
    52 = @{ D=3900; E=10700; F=5900; G=2400 }
    54 = @{ D=110100; E=129700; F=118300; G=58800; H=46100; I=37100 }
    57 = @{ D=8400; E=10400; F=8200; G=9300; H=7500; I=4800 }
    58 = @{ G=3200; H=2700; I=3500 }
    59 = @{ D=9900; E=9800; F=5800; G=2200; H=2500 }
    60 = @{ D=18200; E=20200; F=14000; G=14700; H=12700; I=10100 }
    61 = @{ G=4600; H=7800; I=5000 }
    62 = @{ D=16400; E=10200; F=4600; G=39600; H=32600; I=24900 }
    66 = @{ D=34700; E=30400; F=18600; G=58900; H=53000; I=40000 }
    72 = @{ D=-11000; E=-106000; F=-103000; G=-3100; H=-7900; I=-3000 }
    76 = @{ D=75400; E=99300; F=99800; H=-6900 }
    81 = @{ D=-21200; E=-2700; F=-97800; G=4800; H=-4600; I=-3200 }
    83 = @{ D=7400; E=6600; F=4800; G=4100; H=3500 }
    89 = @{ D=7300; E=8300; F=11600; G=12600; H=10600 }
    91 = @{ D=-7200; E=-15300; F=-12500; G=-9300; H=-6200; I=-9200 }
    94 = @{ D=-9400; E=-16400; F=-13400; G=-9800; H=-6500; I=-10000 }
    100 = @{ D=1300; F=41100; I=13300 }
    102 = @{ E=-6000; F=39200; G=3300; H=5700; I=3400 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}

